$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: apply the same formatting as an existing date cell (column B) to a
# target cell, then assign the date value. Doing it in this order (format
# first, then value) re-uses the workbook's existing "short date" style
# (numFmtId 14) instead of minting a brand-new custom number format.
function Set-DateCell($targetAddr, $sourceAddr, $dateText) {
    $ws.Range($sourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range($targetAddr).Value = $dateText
}

# Row 20
$ws.Range("A20").Value = "JG"
Set-DateCell "B20" "B19" "9/20/2023"
$ws.Range("C20").Value = "ECON110-03"
$ws.Range("D20").Value = 75

# Row 21
$ws.Range("A21").Value = "JG"
Set-DateCell "B21" "B20" "9/20/2023"
$ws.Range("C21").Value = "THEO200-05"
$ws.Range("D21").Value = 85

# Row 22
$ws.Range("A22").Value = "JG"
Set-DateCell "B22" "B21" "9/21/2023"
$ws.Range("C22").Value = "DS160-01"
$ws.Range("D22").Value = 75
$ws.Range("E22").Value = "Matplotlib"

# Row 23
$ws.Range("A23").Value = "JG"
Set-DateCell "B23" "B22" "9/21/2023"
$ws.Range("C23").Value = "MUSE123-05"
$ws.Range("D23").Value = 75

# Row 24 (no initials in column A)
Set-DateCell "B24" "B23" "9/25/2023"
$ws.Range("C24").Value = "MUSE137-01"
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = "Brass ensemble"

# Row 25
$ws.Range("A25").Value = "JG"
Set-DateCell "B25" "B24" "9/25/2023"
$ws.Range("C25").Value = "ECON110-03"
$ws.Range("D25").Value = "?"
$ws.Range("E25").Value = "Test in testing center"

# Row 26
$ws.Range("A26").Value = "JG"
Set-DateCell "B26" "B25" "9/26/2023"
$ws.Range("C26").Value = "DS160-01"
$ws.Range("D26").Value = 75
$ws.Range("E26").Value = "Dataset work with matplotlib"

# Row 27
$ws.Range("A27").Value = "JG"
Set-DateCell "B27" "B26" "9/26/2023"
$ws.Range("C27").Value = "MUSE123-05"
$ws.Range("D27").Value = 75
$ws.Range("E27").Value = "band"

# Row 28 (no initials in column A)
Set-DateCell "B28" "B27" "9/26/2023"
$ws.Range("C28").Value = "MATH205-03"
$ws.Range("D28").Value = 90
$ws.Range("E28").Value = "Stats quiz"

# Update the view: scroll/selection moved to N14 (matches the saved file's
# sheetView selection after the edits were made).
$ws.Range("N14").Select() | Out-Null
